$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('1996 World Champion', ['{W}{U}{B}{R}{G}', 'Summon — Legend', 'Cannot be the target of spells or effects. World Champion has power and toughness equal to the life total of target opponent.', '{0}: Discard your hand to search your library for 1996 World Champion and reveal it to all players. Shuffle your library and put 1996 World Champion on top of it. Use this ability only at the beginning of your upkeep, and only if 1996 World Champion is in your library.'])"
$ws.Range("A3").Value = "('Fraternal Exaltation', ['{U}{U}{U}{U}', 'Sorcery', 'Sneak into your parents’ closet to get a deck. Your new brother is joining the game.'])"
$ws.Range("A4").Value = "('Phoenix Heart', ['{R}{R}{R}{R}', 'Sorcery', 'Koni and Richard continue the game and play all games as partners forever.'])"
$ws.Range("A5").Value = "('Proposal', ['{W}{W}{W}{W}', 'Sorcery', 'Allows Richard to propose marriage to Lily. If the proposal is accepted both players win; mix the cards in play, both libraries, and both graveyards as a shared deck.'])"
$ws.Range("A6").Value = "('Robot Chicken', ['{4}', 'Artifact Creature — Chicken Construct', 'Whenever you cast a spell, put a 0/1 colorless Egg artifact creature token onto the battlefield.', 'Whenever an Egg you control is put into a graveyard from the battlefield, destroy target artifact or creature.', '2/2'])"
$ws.Range("A7").Value = "('Shichifukujin Dragon', ['{6}{R}{R}{R}', 'Summon — Dragon', 'When Shichifukujin Dragon comes into play, put seven +1/+1 counters on it.', '{R}{R}{R}, Sacrifice two +1/+1 counters: Put three +1/+1 counters on Shichifukujin Dragon at end of turn. Play this ability as a sorcery.'])"
$ws.Range("A8").Value = "('Splendid Genesis', ['{G}{G}{G}{G}', 'Sorcery', 'Shuffle all cards in the game together and deal them into three decks. The game continues with a new player.'])"

$rows = $ws.Range("A9:A33").EntireRow
$rows.Delete()
